$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are treated as Text so numeric-looking strings
# (e.g. "1.00", "24.98") keep their exact original formatting/precision.
$changes = @{
    "D2" = "60.410.23"
    "E2" = "  +3.99%  "
    "D3" = "2.429.81"
    "E3" = "  +3.12%  "
    "E4" = "  +0.05%  "
    "D5" = "556.59"
    "E5" = "  +2.15%  "
    "D6" = "139.23"
    "E6" = "  +3.29%  "
    "D7" = "1.00"
    "E7" = "  -0.04%  "
    "D8" = "0.575"
    "E8" = "  +3.33%  "
    "E9" = "  +4.88%  "
    "E10" = "  +3.42%  "
    "E11" = "  +1.55%  "
    "E12" = "  -2.32%  "
    "D13" = "24.98"
    "E13" = "  +4.70%  "
    "D14" = "2.862.05"
    "E14" = "  +3.16%  "
    "D15" = "60.322.01"
    "E15" = "  +3.94%  "
    "E16" = "  +3.89%  "
    "D17" = "2.444.52"
    "E17" = "  +3.13%  "
    "D18" = "11.42"
    "E18" = "  +5.82%  "
    "D19" = "4.43"
    "E19" = "  +2.87%  "
    "D20" = "335.82"
    "E20" = "  +1.27%  "
    "D21" = "6.85"
    "E21" = "  +2.06%  "
    "E22" = "  -0.08%  "
    "D23" = "65.56"
    "E23" = "  +4.58%  "
    "E24" = "  +3.33%  "
    "D25" = "8.57"
    "E25" = "  +0.91%  "
    "E26" = "  +0.06%  "
    "D27" = "1.35"
    "E27" = "  -0.68%  "
    "E28" = "  +6.34%  "
    "E29" = "  +2.26%  "
    "D30" = "6.31"
    "E30" = "  +2.72%  "
    "D31" = "169.83"
    "E31" = "  -0.57%  "
    "D32" = "18.79"
    "E32" = "  +1.72%  "
    "D33" = "1.04"
    "E33" = "  +1.07%  "
    "E35" = "  +6.17%  "
    "E36" = "  -0.37%  "
    "E37" = "  +0.05%  "
    "E38" = "  +0.31%  "
    "D39" = "39.86"
    "E39" = "  +1.04%  "
    "B40" = "Bittensor"
    "C40" = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
    "D40" = "321.88"
    "E40" = "  +10.40%  "
    "B41" = "PolygonEcosystemToken"
    "C41" = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
    "D41" = "0.418"
    "E41" = "  +10.20%  "
    "E42" = "  +1.52%  "
    "D43" = "142.68"
    "E43" = "  -0.66%  "
    "D44" = "0.0527"
    "E44" = "  +3.64%  "
    "D45" = "0.0961"
    "E45" = "  +2.03%  "
    "D46" = "19.73"
    "E46" = "  +3.32%  "
    "D47" = "0.413"
    "E47" = "  +7.71%  "
    "E48" = "  +1.19%  "
    "E49" = "  +1.81%  "
    "D50" = "17.96"
    "E50" = "  +1.93%  "
    "D51" = "11.06"
    "E51" = "  -0.17%  "
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$addr]
}
